# Refresh cryptos list price (D) / volume-1h (E) columns to the latest scrape.
# Values are plain text cells (not numbers), so numeric-looking Price strings
# are written with a leading apostrophe to keep Excel from auto-converting
# them to floating point numbers (matches the original inline-string content).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''62.766.55'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '''2.675.01'
$ws.Range('E3').Value = '  -2.33%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''553.25'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').Value = '''157.78'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '''0.591'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('E10').Value = '  -2.55%  '
$ws.Range('E11').Value = '  -3.29%  '
$ws.Range('D12').Value = '''5.36'
$ws.Range('E12').Value = '  -5.26%  '
$ws.Range('D13').Value = '''3.151.02'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').Value = '''26.42'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').Value = '''62.695.20'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').Value = '''2.681.30'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '''11.86'
$ws.Range('E18').Value = '  -4.10%  '
$ws.Range('D19').Value = '''4.62'
$ws.Range('E19').Value = '  -3.11%  '
$ws.Range('D20').Value = '''344.54'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').Value = '''6.28'
$ws.Range('E21').Value = '  -4.59%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '''0.506'
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('D24').Value = '''63.31'
$ws.Range('E24').Value = '  -1.41%  '
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').Value = '''0.998'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '''8.22'
$ws.Range('E27').Value = '  -1.52%  '
$ws.Range('D28').Value = '''1.43'
$ws.Range('E28').Value = '  +8.97%  '
$ws.Range('D29').Value = '''0.0₃0852'
$ws.Range('E29').Value = '  -5.54%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('E31').Value = '  -0.97%  '
$ws.Range('D32').Value = '''163.24'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = '''4.90'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').Value = '''1.48'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '''19.47'
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Value = '''352.57'
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('D39').Value = '''0.948'
$ws.Range('E39').Value = '  -3.22%  '
$ws.Range('D40').Value = '''6.20'
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').Value = '''38.38'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '''20.89'
$ws.Range('E43').Value = '  -4.03%  '
$ws.Range('D44').Value = '''20.23'
$ws.Range('E44').Value = '  -3.63%  '
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('D46').Value = '''0.0560'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('D47').Value = '''0.999'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('E49').Value = '  -2.73%  '
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('D51').Value = '''128.64'
$ws.Range('E51').Value = '  -4.76%  '
